$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.02508122093216724
$ws.Range("D2").Value = 1.528858312752999
$ws.Range("F2").Value = -0.005335027006694312

# Row 3
$ws.Range("B3").Value = 2.054458368791535
$ws.Range("D3").Value = 125.2321712657749
$ws.Range("F3").Value = 1.577873015505964

# Row 4
$ws.Range("B4").Value = 0.02508122093216724
$ws.Range("D4").Value = 1.528858312752999
$ws.Range("F4").Value = -0.005335027006694312

# Row 5
$ws.Range("B5").Value = 0.2310722898626279
$ws.Range("D5").Value = 1.103891934739909
$ws.Range("F5").Value = -0.09264808326638718

# Row 6
$ws.Range("B6").Value = 1.07281325584755
$ws.Range("D6").Value = 5.125105659861769
$ws.Range("F6").Value = 0.3029987632388947

# Row 7
$ws.Range("B7").Value = 0.3586701185774284
$ws.Range("D7").Value = 1.71345967690549
$ws.Range("F7").Value = 0.3171355525542741

# Row 8
$ws.Range("B8").Value = 2.926769900689424
$ws.Range("D8").Value = 5.748028375001295
$ws.Range("F8").Value = 2.895519011693499

# Row 9
$ws.Range("B9").Value = 27.8520716728664
$ws.Range("D9").Value = 54.70006311069905
$ws.Range("F9").Value = 27.48494852201295

# Row 10
$ws.Range("B10").Value = 2.926769900689424
$ws.Range("D10").Value = 5.748028375001295
$ws.Range("F10").Value = 2.895519011693499

# Row 12
$ws.Range("B12").Value = 5.390526155778248
$ws.Range("D12").Value = 1.347631538944562
$ws.Range("F12").Value = -5.264481861280244

# Row 13
$ws.Range("B13").Value = 6.205103218842413
$ws.Range("D13").Value = 1.551275804710603
$ws.Range("F13").Value = -6.205103193010602
